$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.175.67"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "1.821.80"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.81%  "
$ws.Range("D5").Value = "'312.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.4233"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").Value = "'0.3683"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'0.07238"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'0.8555"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "'20.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "1.830.03"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "'6.698"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "'0.07080"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "'5.293"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "'89.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'0.000008842"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "'15.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").Value = "27.263.25"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").Value = "'5.118"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'10.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "2.055.10"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'1.981"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").Value = "'152.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'2.205"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("D28").Value = "'18.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "'5.225"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").Value = "'116.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("D31").Value = "'0.08840"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "'1.190"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "'0.7475"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("D34").Value = "'2.957"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "'4.436"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "'1.004"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'1.108"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").Value = "'0.01967"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "'7.275"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "'2.862"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'0.1696"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "'0.5029"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "'8.654"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'10.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4743"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'106.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "'0.06388"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "'1.661"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "'1.881"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.17%  "
